# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns for the first data row (1ab66350-...) on the zh-cn and de-de
# localization status sheets, reflecting a newer handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 09:07:17"
$wsZhCn.Range("K2").Value = "2016-09-05 09:07:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 09:07:22"
$wsDeDe.Range("K2").Value = "2016-09-05 09:07:44"

# The Overview sheet's "Latest HO Xliff Generate Date" reflects the most recent
# Correspond Handoff Datetime across languages for that file; since de-de's
# handoff datetime advanced to 09:07:22 (still the latest), update it here too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 09:07:22"
